# Add a new worksheet "Exercícios 2.1 - 10" by duplicating the structurally
# similar "Exercícios 2.1 - 9" sheet (same layout: Max/Min LP model with a
# 2-variable objective + up to 6 constraints), then overwrite the cells that
# differ for this new exercise, add the Solver add-in's hidden defined
# names scoped to the new sheet, and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("Exercícios 2.1 - 9")

# Duplicate the sheet right after "Exercícios 2.1 - 9" -- this carries over
# styles, column widths, merged cells, page setup, etc. "for free".
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Exercícios 2.1 - 10"

# ---- Header / objective function row -------------------------------------
$ws.Range("B3").Value = "Max Z = 0,12x1 +0,35x2"

# B5 keeps the shared "Restrições:" label (index 12) unchanged from the copy.

$ws.Range("B6").Value = "x1 + x2 <= 70000"
$ws.Range("C6").ClearContents()

$ws.Range("B7").Value = "0,4x1 +0,2x2 <= 30000"
$ws.Range("C7").ClearContents()

$ws.Range("B8").Value = "x1 >= 0"
$ws.Range("C8").ClearContents()

$ws.Range("B9").Value = "x2 >= 0"

# ---- Objective coefficients / value ---------------------------------------
$ws.Range("H5").Value = 0.12
$ws.Range("I5").Value = 0.35

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 70000

$ws.Range("H7").Formula = "=(H5*H6)+(I5*I6)"

# ---- Constraints table ------------------------------------------------
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Formula = "=(H11*H6)+(I11*I6)"
$ws.Range("K11").Value = 70000

$ws.Range("H12").Value = 0.4
$ws.Range("I12").Value = 0.2
$ws.Range("J12").Formula = "=(H12*H6)+(I12*I6)"
$ws.Range("K12").Value = 30000

$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("J13").Formula = "=(H13*H6)+(I13*I6)"
$ws.Range("K13").Value = 0

$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Formula = "=(H14*H6)+(I14*I7)"
$ws.Range("K14").Value = 0

# ---- Solver add-in hidden defined names, scoped to the new sheet ----------
function Add-SolverName($sheet, $name, $refersTo) {
    $sheet.Names.Add($name, $refersTo)
    $n = $wb.Names.Item($sheet.Name + "!" + $name)
    $n.Visible = $false
}

Add-SolverName $ws "solver_adj" "='Exercícios 2.1 - 10'!`$H`$6:`$I`$6"
Add-SolverName $ws "solver_cvg" "=0.0001"
Add-SolverName $ws "solver_drv" "=2"
Add-SolverName $ws "solver_eng" "=2"
Add-SolverName $ws "solver_est" "=1"
Add-SolverName $ws "solver_itr" "=2147483647"
Add-SolverName $ws "solver_lhs1" "='Exercícios 2.1 - 10'!`$J`$11"
Add-SolverName $ws "solver_lhs2" "='Exercícios 2.1 - 10'!`$J`$12"
Add-SolverName $ws "solver_lhs3" "='Exercícios 2.1 - 10'!`$J`$13"
Add-SolverName $ws "solver_lhs4" "='Exercícios 2.1 - 10'!`$J`$14"
Add-SolverName $ws "solver_lhs5" "='Exercícios 2.1 - 10'!`$J`$15"
Add-SolverName $ws "solver_lhs6" "='Exercícios 2.1 - 10'!`$J`$16"
Add-SolverName $ws "solver_mip" "=2147483647"
Add-SolverName $ws "solver_mni" "=30"
Add-SolverName $ws "solver_mrt" "=0.075"
Add-SolverName $ws "solver_msl" "=2"
Add-SolverName $ws "solver_neg" "=1"
Add-SolverName $ws "solver_nod" "=2147483647"
Add-SolverName $ws "solver_num" "=4"
Add-SolverName $ws "solver_nwt" "=1"
Add-SolverName $ws "solver_opt" "='Exercícios 2.1 - 10'!`$H`$7"
Add-SolverName $ws "solver_pre" "=0.000001"
Add-SolverName $ws "solver_rbv" "=2"
Add-SolverName $ws "solver_rel1" "=1"
Add-SolverName $ws "solver_rel2" "=1"
Add-SolverName $ws "solver_rel3" "=3"
Add-SolverName $ws "solver_rel4" "=3"
Add-SolverName $ws "solver_rel5" "=1"
Add-SolverName $ws "solver_rel6" "=1"
Add-SolverName $ws "solver_rhs1" "='Exercícios 2.1 - 10'!`$K`$11"
Add-SolverName $ws "solver_rhs2" "='Exercícios 2.1 - 10'!`$K`$12"
Add-SolverName $ws "solver_rhs3" "='Exercícios 2.1 - 10'!`$K`$13"
Add-SolverName $ws "solver_rhs4" "='Exercícios 2.1 - 10'!`$K`$14"
Add-SolverName $ws "solver_rhs5" "='Exercícios 2.1 - 10'!`$K`$15"
Add-SolverName $ws "solver_rhs6" "='Exercícios 2.1 - 10'!`$K`$16"
Add-SolverName $ws "solver_rlx" "=2"
Add-SolverName $ws "solver_rsd" "=0"
Add-SolverName $ws "solver_scl" "=2"
Add-SolverName $ws "solver_sho" "=2"
Add-SolverName $ws "solver_ssz" "=100"
Add-SolverName $ws "solver_tim" "=2147483647"
Add-SolverName $ws "solver_tol" "=0.01"
Add-SolverName $ws "solver_typ" "=1"
Add-SolverName $ws "solver_val" "=0"
Add-SolverName $ws "solver_ver" "=3"

# ---- Make the new sheet the active / selected tab --------------------------
$ws.Activate()
$ws.Range("H7").Select()
